# ---------------------------------------------------------------------------
# Applies the "timeline.docx" edit described by the unified diff:
#   * "Board space ship" -> "Board spaceship" (now struck through)
#   * new paragraph "Scan of planet for best landing site" after
#     "Landing on new planet"
#   * w:lastRenderedPageBreak moved from "Jack runs towards forest" to the
#     preceding underscore-rule paragraph
#   * two struck-through lines get their single run split into three runs
#     with w:proofErr gramStart/gramEnd bookending the word "a"
#   * "Need some exploration of planet to fill space" gains a trailing
#     space run; "Goes on search for more creatures..." is relocated to
#     right after it; "Ancient ruins?" becomes struck through
#   * two new lines ("Chased by spider...", "Captured by spider...") added
#     after the "Tau Ceti" paragraph, plus one extra blank paragraph
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# pkg:package wrapper boilerplate used by Range.InsertXML.
$pkgHeader = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function New-BodyXml([string]$inner) {
    return $pkgHeader + $inner + $pkgFooter
}

# Find a paragraph (1-based Paragraphs index) whose visible text equals
# $searchText exactly (trailing paragraph-mark / cell-mark char is ignored).
function Get-ParaIndexByText($doc, [string]$searchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $searchText) {
            return $i
        }
    }
    return -1
}

# Find a paragraph (1-based Paragraphs index) whose visible text contains
# $needle as a substring. Used for paragraphs built from several runs where
# matching the whole text is more brittle than matching a unique fragment.
function Get-ParaIndexByContains($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# Replace a whole paragraph's contents (and its w:pPr) with freshly built XML.
function Set-ParagraphXml($doc, [int]$index, [string]$innerXml) {
    $para = $doc.Paragraphs.Item($index)
    $para.Range.InsertXML((New-BodyXml $innerXml)) | Out-Null
}

# Insert a brand-new paragraph right after the paragraph currently at
# $index, with body $innerXml; returns the index of the new paragraph.
function Add-ParagraphAfter($doc, [int]$index, [string]$innerXml) {
    $anchor = $doc.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $newIndex = $index + 1
    $newPara = $doc.Paragraphs.Item($newIndex)
    $newPara.Range.InsertXML((New-BodyXml $innerXml)) | Out-Null
    return $newIndex
}

# ===========================================================================
# Work from the bottom of the document upward so earlier text searches keep
# matching paragraphs that have not shifted yet.
# ===========================================================================

# --- Hunk 8: drop the old "Goes on search..." paragraph (it is relocated
#     earlier in the document by Hunk 6). -----------------------------------
$idx = Get-ParaIndexByText $d "Goes on search for more creatures /Captures 2 more for a total of 3 (EXPAND THIS)"
$d.Paragraphs.Item($idx).Range.Delete()

# --- Hunk 7c: one extra blank paragraph amongst the two existing ones,
#     right after the "Tau Ceti..." paragraph's pair of empty paragraphs. --
$idx = Get-ParaIndexByContains $d "Tau Ceti peoples (ancient type)??"
# the two blank paragraphs immediately follow; add a third after them
$blank2 = $idx + 2
Add-ParagraphAfter $d $blank2 "" | Out-Null

# --- Hunk 7a/7b: two new lines after the "Tau Ceti..." paragraph. ---------
$idx = Get-ParaIndexByContains $d "Tau Ceti peoples (ancient type)??"
$chasedXml = '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Chased by spider (killed village)</w:t></w:r></w:p>'
$capturedXml = '<w:p><w:r><w:t>Captured by spider (saved by rock lion)</w:t></w:r></w:p>'
$idx = Add-ParagraphAfter $d $idx $chasedXml
Add-ParagraphAfter $d $idx $capturedXml | Out-Null

# --- Hunk 6c: "Ancient ruins?" becomes struck through. ---------------------
$idx = Get-ParaIndexByText $d "`tAncient ruins?"
Set-ParagraphXml $d $idx '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:tab/><w:t>Ancient ruins?</w:t></w:r></w:p>'

# --- Hunk 6a/6b: trailing space run + relocated "Goes on search..." line. -
$idx = Get-ParaIndexByText $d "Need some exploration of planet to fill space"
Set-ParagraphXml $d $idx '<w:p><w:r><w:t>Need some exploration of planet to fill space</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
Add-ParagraphAfter $d $idx '<w:p><w:r><w:t>Goes on search for more creatures /Captures 2 more for a total of 3 (EXPAND THIS)</w:t></w:r></w:p>' | Out-Null

# --- Hunk 5: split run around "a" with proofErr markers. -------------------
$idx = Get-ParaIndexByContains $d "Jack encounters a creature and attempts to capture it (suspense, breaks a icosahedron)"
$xml5 = '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Jack encounters a creature and attempts to capture it (suspense, breaks </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>a</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> icosahedron)</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d $idx $xml5

# --- Hunk 4: split run around "a" with proofErr markers. -------------------
$idx = Get-ParaIndexByContains $d "Upon leaving forest, jack encounters nick, who uses a icosahedron"
$xml4 = '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Upon leaving forest, jack encounters nick, who uses </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>a</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    "<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=`"preserve`"> icosahedron and unleashes a $([char]0x201C)tamed$([char]0x201D) creature</w:t></w:r>" +
    '</w:p>'
Set-ParagraphXml $d $idx $xml4

# --- Hunk 3b: remove w:lastRenderedPageBreak from "Jack runs towards forest"
$idx = Get-ParaIndexByText $d "Jack runs towards forest"
Set-ParagraphXml $d $idx '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Jack runs towards forest</w:t></w:r></w:p>'

# --- Hunk 3a: add w:lastRenderedPageBreak to the underscore-rule paragraph.
$idx = Get-ParaIndexByText $d "______________________________________________________________________"
Set-ParagraphXml $d $idx '<w:p><w:r><w:lastRenderedPageBreak/><w:t>______________________________________________________________________</w:t></w:r></w:p>'

# --- Hunk 2: new paragraph after "Landing on new planet". ------------------
$idx = Get-ParaIndexByText $d "Landing on new planet"
Add-ParagraphAfter $d $idx '<w:p><w:r><w:tab/><w:t>Scan of planet for best landing site</w:t></w:r></w:p>' | Out-Null

# --- Hunk 1: "Board space ship" -> "Board spaceship" (struck through). -----
$idx = Get-ParaIndexByText $d "Board space ship"
Set-ParagraphXml $d $idx '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Board spaceship</w:t></w:r></w:p>'

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
